# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets,
# matching the regenerated data published at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 3061
$wsExhibit.Range("F4").Value = 108
$wsExhibit.Range("F6").Value = 1820
$wsExhibit.Range("F7").Value = 44
$wsExhibit.Range("F12").Value = 153

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3061
$wsAll.Range("F5").Value = 108
$wsAll.Range("F7").Value = 1820
$wsAll.Range("F8").Value = 44
$wsAll.Range("F13").Value = 153
